$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.474.53"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.66"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.10"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0624"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.98"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.47"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.627.95"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.71"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.469.90"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.23"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +5.47%  "
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.120"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.218.80"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.796"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.763.21"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.84"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.85"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.54"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  -0.75%  "
